$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "AgentTestCases"
$ws.Range("B4").Value = "Agent Related Test Cases"
$ws.Range("C4").Value = "Y"

$ws.Range("C9").Select()
